# Commit: "refactor code and structure files"
#
# Effective data changes in the target diff:
#   - A1 (date serial) updated from 45406 -> 45436 (one month later)
#   - D29..D32 (unit prices) recalculated to new values
#
# (The mergeCells list in the target XML is the same 12 ranges, just
#  re-ordered by the authoring tool on save; the merged regions
#  themselves are unchanged, so no merge/unmerge is required here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436

$ws.Range("D29").Value = 13023.612
$ws.Range("D30").Value = 15370.992
$ws.Range("D31").Value = 18326.952
$ws.Range("D32").Value = 21816.167
